$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top of the data (row 2), pushing all existing
# data rows down by one. Excel extends formatting from the row below into
# the newly inserted row, and the used range grows accordingly (so the
# former last data row duplicates down into the new last row, 95).
$ws.Rows(2).Insert()

# Helper used to create date-look-alike text ("DD-MM-YYYY") without
# Excel's automatic text-to-date conversion kicking in when assigning
# directly via Range.Value. We build the literal text with a formula in a
# scratch cell, then paste-special (values only) it into the destination
# cell so the destination keeps a plain text value/style (no
# reinterpretation as a date, no "number stored as text" quote-prefix
# marker being introduced).
function Set-LiteralText($cellAddr, $text) {
    $ws.Range("ZZ200").Formula = "=""" + $text + """"
    $ws.Range("ZZ200").Copy()
    $ws.Range($cellAddr).PasteSpecial(-4163)
    $ws.Range("ZZ200").ClearContents()
}

# Populate the newly inserted row with the latest day's price entry
# (same circular/price/link as the rest of the 01-11-2025 group, just a
# newer "as of" date).
Set-LiteralText "A2" "08-11-2025"
$ws.Range("B2").Value = "ALUMINIUM INGOT"
$ws.Range("C2").Value = "IE07"
$ws.Range("D2").Value = 297.15
Set-LiteralText "E2" "01-11-2025"
$ws.Range("F2").Value = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-11-2025.pdf"

# The native row-insert operation shifts hyperlink relationship ids (rId)
# inconsistently with respect to the cells' actual text/ref positions
# (off-by-one drift starting a few rows down), and it never creates a
# hyperlink object for the brand-new duplicated last row (95) at all.
# Rather than rely on that shifted mapping, drop every hyperlink in the
# column and recreate one per data row directly from each cell's own
# (already-correct) displayed text, so every link's target matches what
# is shown in the cell.
$ws.Cells.Hyperlinks.Delete()
for ($r = 2; $r -le 95; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $url = $cell.Text
    $ws.Hyperlinks.Add($cell, $url)
}
